# Apply updated crypto price/volume figures (and the WrappedEther/BinanceUSD row-order swap)
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "28.319.65"),
    @("E2", "  +0.21%  "),
    @("D3", "1.863.03"),
    @("E3", "  -1.15%  "),
    @("D4", "1.021"),
    @("E4", "  +1.19%  "),
    @("D5", "315.87"),
    @("E5", "  +0.13%  "),
    @("D6", "1.016"),
    @("E6", "  +0.67%  "),
    @("D7", "0.5094"),
    @("E7", "  -0.92%  "),
    @("D8", "0.3945"),
    @("E8", "  +0.62%  "),
    @("D9", "0.08476"),
    @("E9", "  +0.90%  "),
    @("D10", "1.106"),
    @("E10", "  -1.60%  "),
    @("D11", "41.83"),
    @("E11", "  +0.10%  "),
    @("D12", "6.221"),
    @("E12", "  -0.77%  "),
    @("D13", "20.37"),
    @("E13", "  -1.46%  "),
    @("B14", "BinanceUSD"),
    @("C14", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"),
    @("D14", "1.021"),
    @("E14", "  +1.08%  "),
    @("B15", "WrappedEther"),
    @("C15", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"),
    @("D15", "1.801.88"),
    @("E15", "  -4.77%  "),
    @("D16", "7.179"),
    @("E16", "  -1.47%  "),
    @("D17", "0.00001114"),
    @("E17", "  +0.74%  "),
    @("D18", "90.38"),
    @("E18", "  -1.05%  "),
    @("D19", "0.06727"),
    @("E19", "  +0.18%  "),
    @("D20", "17.64"),
    @("E20", "  -1.21%  "),
    @("D21", "1.015"),
    @("E21", "  +0.52%  "),
    @("D22", "5.939"),
    @("E22", "  -1.87%  "),
    @("D23", "28.312.57"),
    @("E23", "  +0.06%  "),
    @("D24", "11.10"),
    @("E24", "  -0.61%  "),
    @("D25", "2.275"),
    @("E25", "  +0.66%  "),
    @("D26", "161.49"),
    @("E26", "  +1.13%  "),
    @("D27", "2.021.45"),
    @("E27", "  -3.56%  "),
    @("D28", "20.62"),
    @("E28", "  -0.95%  "),
    @("D29", "2.345"),
    @("E29", "  -5.03%  "),
    @("D30", "127.07"),
    @("E30", "  +0.27%  "),
    @("D31", "0.1046"),
    @("E31", "  -0.39%  "),
    @("E32", "  -0.61%  "),
    @("D33", "5.767"),
    @("E33", "  -1.97%  "),
    @("D34", "3.625"),
    @("E34", "  -0.14%  "),
    @("D35", "0.02415"),
    @("E35", "  -1.69%  "),
    @("D36", "0.06451"),
    @("E36", "  -1.97%  "),
    @("D37", "0.2177"),
    @("E37", "  -1.82%  "),
    @("D38", "8.816"),
    @("E38", "  -8.26%  "),
    @("D39", "1.256"),
    @("E39", "  +1.36%  "),
    @("E40", "  -1.99%  "),
    @("D41", "0.6349"),
    @("E41", "  -2.37%  "),
    @("D42", "4.980"),
    @("E42", "  -0.82%  "),
    @("E43", "  -1.12%  "),
    @("D44", "0.6001"),
    @("E44", "  -1.67%  "),
    @("D45", "12.95"),
    @("E45", "  -1.31%  "),
    @("D46", "3.697"),
    @("E46", "  +0.04%  "),
    @("D47", "1.211"),
    @("E47", "  -5.57%  "),
    @("D48", "1.983"),
    @("E48", "  -1.70%  "),
    @("D49", "1.201"),
    @("E49", "  -2.82%  "),
    @("D50", "120.57"),
    @("E50", "  -0.66%  "),
    @("D51", "0.06822"),
    @("E51", "  -1.51%  ")
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    # Leading apostrophe forces Excel's literal-text entry so numeric-looking
    # strings (e.g. "1.021") are not auto-converted to numbers, matching the
    # inline-string cells in the source file; resetting the style afterward
    # avoids leaving a stray Text-format style behind.
    $ws.Range($cellRef).Value = "'" + $newValue
    $ws.Range($cellRef).Style = "Normal"
}
